# FEATURE: Fecth Busince Confidence file from INE
#
# Rename the "Task Code" column header to "Job Code", reword the job's
# description/error text for the new Business Confidence (INE) fetch job,
# and refresh the last-run date/time stamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column header B3: "Task Code" -> "Job Code"
$ws.Range("B3").Value = "Job Code"

# D4: "Economic Activity Index update" -> "Economic Activity update"
$ws.Range("D4").Value = "Economic Activity update"

# F4: old format-error message -> new DB/date failure message
$ws.Range("F4").Value = "Database failed to get business-confidence-aggregate update date"

# G4: last-run timestamp (serial date) refreshed
$ws.Range("G4").Value = 44831.60458203166
